# Update Handback status report timestamps / priority as part of regenerating
# the Handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# 30960fdb... file, used both in row 2 and row 5 (same timestamp value).
$wsOverview.Range("G2").Value = "2016-09-03 16:21:14"
$wsOverview.Range("G5").Value = "2016-09-03 16:21:14"

# zh-cn sheet: Priority changed from "ht" to "mt" for rows 2 and 5.
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H) and Correspond Handback
# DateTime (K) refreshed with new timestamps.
$wsZhCn.Range("H2").Value = "2016-09-03 16:21:08"
$wsZhCn.Range("H5").Value = "2016-09-03 16:21:08"
$wsZhCn.Range("K2").Value = "2016-09-03 16:21:25"
$wsZhCn.Range("K5").Value = "2016-09-03 16:21:25"

# de-de sheet: Priority changed from "ht" to "mt" for rows 2 and 5.
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# de-de sheet: Correspond Handback DateTime (K) refreshed with new timestamp.
$wsDeDe.Range("K2").Value = "2016-09-03 16:21:32"
$wsDeDe.Range("K5").Value = "2016-09-03 16:21:32"
